$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Resolving-Mac" sending cluster (previously rows 14-17) was removed from
# this LR-pair sheet as part of a TPM recompute. Deleting the rows also lets Excel
# prune the now-orphaned "Resolving-Mac" shared string and re-pack the shared
# string table / type="s" indices exactly like the authoritative diff shows.
$ws.Rows("14:17").Delete()

# Recalculated TPM-derived numeric values for the surviving rows (2-13); only
# cells whose value actually changed are listed, keyed by row -> column -> value.
$updates = @{
    2 = @{ G=0.1004046666666667; H=0.301214; I=0.3697907746891698; J=0.44646652694238; M=11.839004; N=35.517012; O=0.4816941403820247; P=0.5139419866672059; Q=1.188691250285333; R=10.698221252568; S=0.1781260493351026; T=0.2294578938371744 }
    3 = @{ G=0.1004046666666667; H=0.301214; I=0.3697907746891698; J=0.44646652694238; O=0.3280382642169655; P=0.3499993524538634; Q=0.8095099810104445; R=7.285589829094; S=0.1213055238524822; T=0.1562629953221584 }
    4 = @{ G=0.1004046666666667; H=0.301214; I=0.3697907746891698; J=0.44646652694238; K=2; L=0.6666666666666666; M=0.049877; N=0.149631; O=0.002029347962027401; P=0.002165206166751885; Q=0.005007883559333333; R=0.045070952034; S=0.0007504341549920005; T=0.0009666920773839379 }
    5 = @{ G=0.1004046666666667; H=0.301214; I=0.3697907746891698; J=0.44646652694238; M=4.6264905; N=9.252981; O=0.1882382474389825; P=0.1338934547121788; Q=0.464521236489; R=2.787127418934; S=0.06960876734659298; T=0.0597789457056633 }
    6 = @{ I=0.1149919741207596; J=0.1388354464902425; M=11.839004; N=35.517012; O=0.4816941403820247; P=0.5139419866672059; Q=0.3696413292226667; R=3.326771963004; S=0.05539096012493132; T=0.0713533651890238 }
    7 = @{ I=0.1149919741207596; J=0.1388354464902425; O=0.3280382642169655; P=0.3499993524538634; S=0.03772176758945619; T=0.04859231636922789 }
    8 = @{ I=0.1149919741207596; J=0.1388354464902425; K=2; L=0.6666666666666666; M=0.049877; N=0.149631; O=0.002029347962027401; P=0.002165206166751885; Q=0.001557276319666667; R=0.014015486877; S=0.0002333587283314711; T=0.0003006073649044245 }
    9 = @{ I=0.1149919741207596; J=0.1388354464902425; M=4.6264905; N=9.252981; O=0.1882382474389825; P=0.1338934547121788; Q=0.1444498285545; R=0.866698971327; S=0.02164588767804062; T=0.0185891575670864 }
    10 = @{ G=0.1398905; H=0.279781; I=0.5152172511900707; J=0.4146980265673775; M=11.839004; N=35.517012; O=0.4816941403820247; P=0.5139419866672059; Q=1.656164189062; R=9.936985134372; S=0.2481771309219908; T=0.2131307276410077 }
    11 = @{ G=0.1398905; H=0.279781; I=0.5152172511900707; J=0.4146980265673775; O=0.3280382642169655; P=0.3499993524538634; Q=1.127863472466833; R=6.767180834801001; S=0.1690109727750271; T=0.1451440407624772 }
    12 = @{ G=0.1398905; H=0.279781; I=0.5152172511900707; J=0.4146980265673775; K=2; L=0.6666666666666666; M=0.049877; N=0.149631; O=0.002029347962027401; P=0.002165206166751885; Q=0.006977318468500001; R=0.04186391081100001; S=0.001045555078703929; T=0.0008979067244635228 }
    13 = @{ G=0.1398905; H=0.279781; I=0.5152172511900707; J=0.4146980265673775; M=4.6264905; N=9.252981; O=0.1882382474389825; P=0.1338934547121788; Q=0.64720206929025; R=2.588808277161; S=0.09698359241434894; T=0.05552535143942906 }
}

foreach ($rowNum in $updates.Keys) {
    $rowValues = $updates[$rowNum]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$rowNum").Value = $rowValues[$col]
    }
}
